$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-fte-count"
$ws1.Range("B3").Value = "8.0.0"
$ws1.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$ws1.Range("B9").Value = "LinuxForHealth Team"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("AI2").Value = ""
$ws2.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-fte-count"
